$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

$ws.Range("D2").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Ramírez González" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 67890123 - 1681735068" + $nl + "Roles: ['Solicitante', 'Gestor 1', 'Recepción', 'Gestor 2', 'Administrador']"

$ws.Range("D3").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Ramírez González" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 67890123 - 1681735068" + $nl + "Roles: ['Solicitante', 'Gestor 1', 'Recepción', 'Gestor 2', 'Administrador']"
$ws.Range("F3").Value = "NO : no se encontraron resultados en Usuarios para ['gamerboy42 - 1681735068', 'Solicitante Gestor 1 Recepción Gestor 2 Administrador']"
$ws.Range("G3").Value = "FAILED"

$ws.Range("D4").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Ramírez González" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 67890123 - 1681735068" + $nl + "Roles: ['Solicitante', 'Gestor 1', 'Recepción', 'Gestor 2', 'Administrador']"
$ws.Range("F4").Value = "SI : ['gamerboy42 - 1681735068', 'Martín Andrés', 'Ramírez González', 'Pasaporte', '67890123 - 1681735068', True, ['Administrador', 'Gestor 1', 'Gestor 2', 'Recepción', 'Solicitante']] y ['gamerboy42 - 1681735068', 'Martín Andrés', 'Ramírez González', 'Pasaporte', '67890123 - 1681735068', True, ['Administrador', 'Gestor 1', 'Gestor 2', 'Recepción', 'Solicitante']] coinciden"

$ws.Range("D5").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: FG789012 - 1681735110" + $nl + "Roles: ['Gestor 2', 'Recepción']"

$ws.Range("D6").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: FG789012 - 1681735110" + $nl + "Roles: ['Gestor 2', 'Recepción']"
$ws.Range("F6").Value = "NO : no se encontraron resultados en Usuarios para ['mysticalunicorn88 - 1681735110', 'Gestor 2 Recepción']"
$ws.Range("G6").Value = "FAILED"

$ws.Range("D7").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: FG789012 - 1681735110" + $nl + "Roles: ['Gestor 2', 'Recepción']"
$ws.Range("F7").Value = "SI : ['mysticalunicorn88 - 1681735110', 'Ana Isabel', 'Martínez Sánchez', 'Cédula de Extranjería', 'FG789012 - 1681735110', True, ['Gestor 2', 'Recepción']] y ['mysticalunicorn88 - 1681735110', 'Ana Isabel', 'Martínez Sánchez', 'Cédula de Extranjería', 'FG789012 - 1681735110', True, ['Gestor 2', 'Recepción']] coinciden"

$ws.Range("D8").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: 1234567A - 1681735152" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Administrador', 'Gestor 2', 'Solicitante']"

$ws.Range("D9").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: 1234567A - 1681735152" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Administrador', 'Gestor 2', 'Solicitante']"
$ws.Range("F9").Value = "SI : se encontró un resultado en Usuarios que coincide con [['lovetoswim99 - 1681735152', 'Recepción Gestor 1 Administrador Gestor 2 Solicitante']]"

$ws.Range("D10").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: 1234567A - 1681735152" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Administrador', 'Gestor 2', 'Solicitante']"
$ws.Range("F10").Value = "SI : ['lovetoswim99 - 1681735152', 'Martín Andrés', 'Martínez Sánchez', 'T.I.', '1234567A - 1681735152', True, ['Administrador', 'Gestor 1', 'Gestor 2', 'Recepción', 'Solicitante']] y ['lovetoswim99 - 1681735152', 'Martín Andrés', 'Martínez Sánchez', 'T.I.', '1234567A - 1681735152', True, ['Administrador', 'Gestor 1', 'Gestor 2', 'Recepción', 'Solicitante']] coinciden"

$ws.Range("D11").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Gómez Rodríguez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 4567890B - 1681735192" + $nl + "Roles: ['Solicitante', 'Gestor 2', 'Administrador', 'Recepción']"

$ws.Range("D12").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Gómez Rodríguez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 4567890B - 1681735192" + $nl + "Roles: ['Solicitante', 'Gestor 2', 'Administrador', 'Recepción']"
$ws.Range("F12").Value = "SI : se encontró un resultado en Usuarios que coincide con [['hikingfanatic33 - 1681735192', 'Solicitante Gestor 2 Administrador Recepción']]"

$ws.Range("D13").Value = "Nombres: Martín Andrés" + $nl + "Apellidos: Gómez Rodríguez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 4567890B - 1681735192" + $nl + "Roles: ['Solicitante', 'Gestor 2', 'Administrador', 'Recepción']"
$ws.Range("F13").Value = "SI : ['hikingfanatic33 - 1681735192', 'Martín Andrés', 'Gómez Rodríguez', 'C.C.', '4567890B - 1681735192', True, ['Administrador', 'Gestor 2', 'Recepción', 'Solicitante']] y ['hikingfanatic33 - 1681735192', 'Martín Andrés', 'Gómez Rodríguez', 'C.C.', '4567890B - 1681735192', True, ['Administrador', 'Gestor 2', 'Recepción', 'Solicitante']] coinciden"

$ws.Range("D14").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: García Pérez" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 8901234C - 1681735237" + $nl + "Roles: ['Recepción']"

$ws.Range("D15").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: García Pérez" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 8901234C - 1681735237" + $nl + "Roles: ['Recepción']"
$ws.Range("F15").Value = "SI : se encontró un resultado en Usuarios que coincide con [['teadrinker12 - 1681735237', 'Recepción']]"

$ws.Range("D16").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: García Pérez" + $nl + "Tipo Doc: Pasaporte" + $nl + "Num Doc: 8901234C - 1681735237" + $nl + "Roles: ['Recepción']"
$ws.Range("F16").Value = "SI : ['teadrinker12 - 1681735237', 'Ana Isabel', 'García Pérez', 'Pasaporte', '8901234C - 1681735237', True, ['Recepción']] y ['teadrinker12 - 1681735237', 'Ana Isabel', 'García Pérez', 'Pasaporte', '8901234C - 1681735237', True, ['Recepción']] coinciden"

$ws.Range("D17").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 2345678D - 1681735278" + $nl + "Roles: ['Solicitante']"

$ws.Range("D18").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 2345678D - 1681735278" + $nl + "Roles: ['Solicitante']"
$ws.Range("F18").Value = "SI : se encontró un resultado en Usuarios que coincide con [['familyman77 - 1681735278', 'Solicitante']]"

$ws.Range("D19").Value = "Nombres: Ana Isabel" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 2345678D - 1681735278" + $nl + "Roles: ['Solicitante']"
$ws.Range("F19").Value = "SI : ['familyman77 - 1681735278', 'Ana Isabel', 'Martínez Sánchez', 'Cédula de Extranjería', '2345678D - 1681735278', True, ['Solicitante']] y ['familyman77 - 1681735278', 'Ana Isabel', 'Martínez Sánchez', 'Cédula de Extranjería', '2345678D - 1681735278', True, ['Solicitante']] coinciden"

$ws.Range("D20").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Castro Ruiz" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 5678901E - 1681735319" + $nl + "Roles: ['Administrador', 'Recepción', 'Gestor 2']"

$ws.Range("D21").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Castro Ruiz" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 5678901E - 1681735319" + $nl + "Roles: ['Administrador', 'Recepción', 'Gestor 2']"
$ws.Range("F21").Value = "NO : no se encontraron resultados en Usuarios para ['doglover88 - 1681735319', 'Administrador Recepción Gestor 2']"
$ws.Range("G21").Value = "FAILED"

$ws.Range("D22").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Castro Ruiz" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: 5678901E - 1681735319" + $nl + "Roles: ['Administrador', 'Recepción', 'Gestor 2']"
$ws.Range("F22").Value = "SI : ['doglover88 - 1681735319', 'Santiago Alejandro', 'Castro Ruiz', 'C.C.', '5678901E - 1681735319', True, ['Administrador', 'Gestor 2', 'Recepción']] y ['doglover88 - 1681735319', 'Santiago Alejandro', 'Castro Ruiz', 'C.C.', '5678901E - 1681735319', True, ['Administrador', 'Gestor 2', 'Recepción']] coinciden"

$ws.Range("D23").Value = "Nombres: Valentina Victoria" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 9012345F - 1681735360" + $nl + "Roles: ['Gestor 1', 'Recepción', 'Solicitante']"

$ws.Range("D24").Value = "Nombres: Valentina Victoria" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 9012345F - 1681735360" + $nl + "Roles: ['Gestor 1', 'Recepción', 'Solicitante']"
$ws.Range("F24").Value = "SI : se encontró un resultado en Usuarios que coincide con [['catlady44 - 1681735360', 'Gestor 1 Recepción Solicitante']]"

$ws.Range("D25").Value = "Nombres: Valentina Victoria" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: Cédula de Extranjería" + $nl + "Num Doc: 9012345F - 1681735360" + $nl + "Roles: ['Gestor 1', 'Recepción', 'Solicitante']"
$ws.Range("F25").Value = "SI : ['catlady44 - 1681735360', 'Valentina Victoria', 'Martínez Sánchez', 'Cédula de Extranjería', '9012345F - 1681735360', True, ['Gestor 1', 'Recepción', 'Solicitante']] y ['catlady44 - 1681735360', 'Valentina Victoria', 'Martínez Sánchez', 'Cédula de Extranjería', '9012345F - 1681735360', True, ['Gestor 1', 'Recepción', 'Solicitante']] coinciden"

$ws.Range("D26").Value = "Nombres: Sofía Elena" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: CD567890 - 1681735403" + $nl + "Roles: ['Solicitante', 'Recepción']"

$ws.Range("D27").Value = "Nombres: Sofía Elena" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: CD567890 - 1681735403" + $nl + "Roles: ['Solicitante', 'Recepción']"
$ws.Range("F27").Value = "NO : no se encontraron resultados en Usuarios para ['birdwatcher22 - 1681735403', 'Solicitante Recepción']"
$ws.Range("G27").Value = "FAILED"

$ws.Range("D28").Value = "Nombres: Sofía Elena" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: T.I." + $nl + "Num Doc: CD567890 - 1681735403" + $nl + "Roles: ['Solicitante', 'Recepción']"
$ws.Range("E28").Value = "EXCEPTION"
$ws.Range("F28").Value = "list index out of range"
$ws.Range("G28").Value = "EXCEPTION"

$ws.Range("D29").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: EF789012 - 1681735443" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Gestor 2']"

$ws.Range("D30").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: EF789012 - 1681735443" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Gestor 2']"
$ws.Range("F30").Value = "SI : se encontró un resultado en Usuarios que coincide con [['butterflykisses44 - 1681735443', 'Recepción Gestor 1 Gestor 2']]"

$ws.Range("D31").Value = "Nombres: Santiago Alejandro" + $nl + "Apellidos: Martínez Sánchez" + $nl + "Tipo Doc: C.C." + $nl + "Num Doc: EF789012 - 1681735443" + $nl + "Roles: ['Recepción', 'Gestor 1', 'Gestor 2']"
$ws.Range("F31").Value = "SI : ['butterflykisses44 - 1681735443', 'Santiago Alejandro', 'Martínez Sánchez', 'C.C.', 'EF789012 - 1681735443', True, ['Gestor 1', 'Gestor 2', 'Recepción']] y ['butterflykisses44 - 1681735443', 'Santiago Alejandro', 'Martínez Sánchez', 'C.C.', 'EF789012 - 1681735443', True, ['Gestor 1', 'Gestor 2', 'Recepción']] coinciden"
